$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 20,10

$arr[0,0] = -18.28040377994633
$arr[0,1] = 1.707305730903584
$arr[0,2] = -18.28040377994633
$arr[0,3] = -18.28040377994633
$arr[0,4] = -18.28040377994633
$arr[0,5] = -18.28040377994633
$arr[0,6] = -18.28040377994633
$arr[0,7] = -18.28040377994633
$arr[0,8] = -18.28040377994633
$arr[0,9] = -18.28040377994633

$arr[1,0] = -18.28040377994633
$arr[1,1] = -18.28040377994633
$arr[1,2] = -18.28040377994633
$arr[1,3] = -18.28040377994633
$arr[1,4] = -18.28040377994633
$arr[1,5] = -18.28040377994633
$arr[1,6] = -18.28040377994633
$arr[1,7] = 2.701604494009925
$arr[1,8] = -18.28040377994633
$arr[1,9] = -18.28040377994633

$arr[2,0] = -18.28040377994633
$arr[2,1] = 1.88095079992735
$arr[2,2] = 2.42509446212585
$arr[2,3] = -18.28040377994633
$arr[2,4] = 3.593803525350883
$arr[2,5] = -18.28040377994633
$arr[2,6] = 1.725479771945166
$arr[2,7] = -18.28040377994633
$arr[2,8] = 2.184808866455103
$arr[2,9] = -18.28040377994633

$arr[3,0] = -18.28040377994633
$arr[3,1] = 2.186780835123331
$arr[3,2] = -18.28040377994633
$arr[3,3] = -18.28040377994633
$arr[3,4] = -18.28040377994633
$arr[3,5] = 2.809295610179022
$arr[3,6] = -18.28040377994633
$arr[3,7] = -18.28040377994633
$arr[3,8] = -18.28040377994633
$arr[3,9] = -18.28040377994633

$arr[4,0] = -18.28040377994633
$arr[4,1] = -18.28040377994633
$arr[4,2] = -18.28040377994633
$arr[4,3] = -18.28040377994633
$arr[4,4] = -18.28040377994633
$arr[4,5] = -18.28040377994633
$arr[4,6] = -18.28040377994633
$arr[4,7] = -18.28040377994633
$arr[4,8] = -18.28040377994633
$arr[4,9] = -18.28040377994633

$arr[5,0] = 2.575006905684166
$arr[5,1] = -18.28040377994633
$arr[5,2] = -18.28040377994633
$arr[5,3] = -18.28040377994633
$arr[5,4] = -18.28040377994633
$arr[5,5] = -18.28040377994633
$arr[5,6] = -18.28040377994633
$arr[5,7] = -18.28040377994633
$arr[5,8] = -18.28040377994633
$arr[5,9] = -18.28040377994633

$arr[6,0] = -18.28040377994633
$arr[6,1] = -18.28040377994633
$arr[6,2] = -18.28040377994633
$arr[6,3] = -18.28040377994633
$arr[6,4] = -18.28040377994633
$arr[6,5] = -18.28040377994633
$arr[6,6] = -18.28040377994633
$arr[6,7] = -18.28040377994633
$arr[6,8] = -18.28040377994633
$arr[6,9] = -18.28040377994633

$arr[7,0] = 3.811594855457892
$arr[7,1] = -18.28040377994633
$arr[7,2] = -18.28040377994633
$arr[7,3] = -18.28040377994633
$arr[7,4] = -18.28040377994633
$arr[7,5] = -18.28040377994633
$arr[7,6] = -18.28040377994633
$arr[7,7] = -18.28040377994633
$arr[7,8] = -18.28040377994633
$arr[7,9] = -18.28040377994633

$arr[8,0] = -18.28040377994633
$arr[8,1] = -18.28040377994633
$arr[8,2] = -18.28040377994633
$arr[8,3] = -18.28040377994633
$arr[8,4] = -18.28040377994633
$arr[8,5] = -18.28040377994633
$arr[8,6] = -18.28040377994633
$arr[8,7] = 1.256679983605928
$arr[8,8] = -18.28040377994633
$arr[8,9] = 1.941401618898104

$arr[9,0] = -18.28040377994633
$arr[9,1] = -18.28040377994633
$arr[9,2] = -18.28040377994633
$arr[9,3] = 4.32192379012758
$arr[9,4] = -18.28040377994633
$arr[9,5] = 2.889291893068965
$arr[9,6] = -18.28040377994633
$arr[9,7] = -18.28040377994633
$arr[9,8] = -18.28040377994633
$arr[9,9] = 1.866260599637634

$arr[10,0] = -18.28040377994633
$arr[10,1] = -18.28040377994633
$arr[10,2] = -18.28040377994633
$arr[10,3] = -18.28040377994633
$arr[10,4] = -18.28040377994633
$arr[10,5] = -18.28040377994633
$arr[10,6] = -18.28040377994633
$arr[10,7] = -18.28040377994633
$arr[10,8] = -18.28040377994633
$arr[10,9] = -18.28040377994633

$arr[11,0] = -18.28040377994633
$arr[11,1] = -18.28040377994633
$arr[11,2] = -18.28040377994633
$arr[11,3] = -18.28040377994633
$arr[11,4] = -18.28040377994633
$arr[11,5] = -18.28040377994633
$arr[11,6] = -18.28040377994633
$arr[11,7] = -18.28040377994633
$arr[11,8] = 1.998899602683313
$arr[11,9] = 1.960737983573392

$arr[12,0] = -18.28040377994633
$arr[12,1] = -18.28040377994633
$arr[12,2] = 1.254858305274797
$arr[12,3] = -18.28040377994633
$arr[12,4] = -18.28040377994633
$arr[12,5] = -18.28040377994633
$arr[12,6] = -18.28040377994633
$arr[12,7] = -18.28040377994633
$arr[12,8] = -18.28040377994633
$arr[12,9] = 2.16562643757046

$arr[13,0] = -18.28040377994633
$arr[13,1] = -18.28040377994633
$arr[13,2] = 0.8590601391902465
$arr[13,3] = -18.28040377994633
$arr[13,4] = -18.28040377994633
$arr[13,5] = -18.28040377994633
$arr[13,6] = -18.28040377994633
$arr[13,7] = -18.28040377994633
$arr[13,8] = -18.28040377994633
$arr[13,9] = -18.28040377994633

$arr[14,0] = -18.28040377994633
$arr[14,1] = -18.28040377994633
$arr[14,2] = -18.28040377994633
$arr[14,3] = -18.28040377994633
$arr[14,4] = -18.28040377994633
$arr[14,5] = -18.28040377994633
$arr[14,6] = -18.28040377994633
$arr[14,7] = -18.28040377994633
$arr[14,8] = 2.11068139336528
$arr[14,9] = -18.28040377994633

$arr[15,0] = -18.28040377994633
$arr[15,1] = 1.93254969214165
$arr[15,2] = 2.195275233910112
$arr[15,3] = -18.28040377994633
$arr[15,4] = -18.28040377994633
$arr[15,5] = -18.28040377994633
$arr[15,6] = 1.368025780525577
$arr[15,7] = 2.002590719948174
$arr[15,8] = 2.092946763341547
$arr[15,9] = -18.28040377994633

$arr[16,0] = -18.28040377994633
$arr[16,1] = -18.28040377994633
$arr[16,2] = -18.28040377994633
$arr[16,3] = -18.28040377994633
$arr[16,4] = -18.28040377994633
$arr[16,5] = -18.28040377994633
$arr[16,6] = 1.612182021411043
$arr[16,7] = 1.257741568013801
$arr[16,8] = 1.521653402784645
$arr[16,9] = -18.28040377994633

$arr[17,0] = -18.28040377994633
$arr[17,1] = -18.28040377994633
$arr[17,2] = 1.726548180314418
$arr[17,3] = -18.28040377994633
$arr[17,4] = -18.28040377994633
$arr[17,5] = -18.28040377994633
$arr[17,6] = 1.534912275764599
$arr[17,7] = 1.528549400312163
$arr[17,8] = -18.28040377994633
$arr[17,9] = -18.28040377994633

$arr[18,0] = -18.28040377994633
$arr[18,1] = 1.103192944182804
$arr[18,2] = 1.344771729510155
$arr[18,3] = -18.28040377994633
$arr[18,4] = 2.986625749536616
$arr[18,5] = -18.28040377994633
$arr[18,6] = 1.918330203235204
$arr[18,7] = 0.8654413184720742
$arr[18,8] = -18.28040377994633
$arr[18,9] = 2.047566534405029

$arr[19,0] = -18.28040377994633
$arr[19,1] = 1.340734122383554
$arr[19,2] = -18.28040377994633
$arr[19,3] = -18.28040377994633
$arr[19,4] = -18.28040377994633
$arr[19,5] = 2.480641757077756
$arr[19,6] = 2.130003058071056
$arr[19,7] = -18.28040377994633
$arr[19,8] = -18.28040377994633
$arr[19,9] = -18.28040377994633

$ws.Range("B2:K21").Value = $arr
